# Updated cryptos list on Fri May 17 16:45:33 UTC 2024 with GitHub Actions
# Refresh the live crypto price/volume snapshot in Sheet1 (rows 2-51).
# Row 40-43 also got re-sorted (coin rank reshuffled) so Coin/Link/Price/Volume
# cells are rewritten together for those rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.205.81'
$ws.Range('E2').Value = '  +2.80%  '

# Row 3
$ws.Range('D3').Value = '3.104.60'
$ws.Range('E3').Value = '  +5.02%  '

# Row 4
$ws.Range('E4').Value = '  +0.09%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Formula = '582.98'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.91%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Formula = '170.35'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +7.54%  '

# Row 8
$ws.Range('D8').Value = '3.098.73'
$ws.Range('E8').Value = '  +4.88%  '

# Row 9
$ws.Range('E9').Value = '  +1.22%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Formula = '6.69'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.44%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Formula = '0.155'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.93%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Formula = '0.483'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +5.64%  '

# Row 13
$ws.Range('E13').Value = '  +2.89%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Formula = '37.05'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +8.84%  '

# Row 15
$ws.Range('E15').Value = '  -0.46%  '

# Row 16
$ws.Range('D16').Value = '3.622.42'
$ws.Range('E16').Value = '  +5.14%  '

# Row 17
$ws.Range('D17').Value = '67.193.74'
$ws.Range('E17').Value = '  +2.64%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Formula = '7.26'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +4.18%  '

# Row 19
$ws.Range('D19').Value = '3.104.50'
$ws.Range('E19').Value = '  +5.03%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Formula = '16.30'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +17.69%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Formula = '472.28'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +5.46%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Formula = '0.717'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +5.52%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Formula = '7.56'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +4.84%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Formula = '84.04'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.40%  '

# Row 25
$ws.Range('E25').Value = '  +9.08%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Formula = '12.95'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +7.29%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Formula = '10.26'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +3.33%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Formula = '8.11'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.59%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Formula = '2.44'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +4.63%  '

# Row 31
$ws.Range('E31').Value = '  +4.84%  '

# Row 32
$ws.Range('E32').Value = '  +5.01%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Formula = '28.49'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +4.04%  '

# Row 34
$ws.Range('E34').Value = '  +5.66%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Formula = '0.999'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.07%  '

# Row 36
$ws.Range('E36').Value = '  +4.28%  '

# Row 37
$ws.Range('E37').Value = '  +3.68%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Formula = '47.37'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +10.74%  '

# Row 39
$ws.Range('E39').Value = '  +6.73%  '

# Row 40
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Formula = '0.319'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +7.01%  '

# Row 41
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Formula = '50.43'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.68%  '

# Row 42
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Formula = '2.94'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +6.14%  '

# Row 43
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Formula = '0.124'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +4.50%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Formula = '8.76'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +3.85%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Formula = '396.93'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.13%  '

# Row 46
$ws.Range('E46').Value = '  +3.61%  '

# Row 47
$ws.Range('D47').Value = '2.772.91'
$ws.Range('E47').Value = '  +1.39%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Formula = '135.31'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +3.59%  '

# Row 49
$ws.Range('E49').Value = '  +0.02%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Formula = '24.88'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +7.39%  '

# Row 51
$ws.Range('E51').Value = '  +5.22%  '
